$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A185:K185").Copy($ws.Range("A186:K186"))
$ws.Range("A185:K185").Copy($ws.Range("A187:K187"))
$ws.Range("A185:K185").Copy($ws.Range("A188:K188"))
$ws.Range("A185:K185").Copy($ws.Range("A189:K189"))
$ws.Range("A185:K185").Copy($ws.Range("A190:K190"))
$ws.Range("A185:K185").Copy($ws.Range("A191:K191"))
$ws.Range("A185:K185").Copy($ws.Range("A192:K192"))
$ws.Range("A185:K185").Copy($ws.Range("A193:K193"))
$ws.Range("A185:K185").Copy($ws.Range("A194:K194"))
$ws.Range("A185:K185").Copy($ws.Range("A195:K195"))
$ws.Range("A185:K185").Copy($ws.Range("A196:K196"))
$ws.Range("A185:K185").Copy($ws.Range("A197:K197"))
$ws.Range("A185:K185").Copy($ws.Range("A198:K198"))
$ws.Range("A185:K185").Copy($ws.Range("A199:K199"))
$ws.Range("A185:K185").Copy($ws.Range("A200:K200"))
$ws.Range("A185:K185").Copy($ws.Range("A201:K201"))
$ws.Range("A185:K185").Copy($ws.Range("A202:K202"))
$ws.Range("A185:K185").Copy($ws.Range("A203:K203"))

$ws.Cells.Item(186,1).Value = 46051.748093703703
$ws.Cells.Item(186,2).Value = 'MH001GMB'
$ws.Cells.Item(186,3).Value = 'Gombe'
$ws.Cells.Item(186,4).Value = 'Kumo'
$ws.Cells.Item(186,5).Value = 'Cowpea White'
$ws.Cells.Item(186,6).Value = 55000
$ws.Cells.Item(186,7).Value = 80
$ws.Cells.Item(186,8).Value = 687.5
$ws.Cells.Item(186,9).Value = 'medium'
$ws.Cells.Item(186,10).Value = 'New'
$ws.Cells.Item(186,11).Value = 690.5

$ws.Cells.Item(187,1).Value = 46051.748703333331
$ws.Cells.Item(187,2).Value = 'MH001GMB'
$ws.Cells.Item(187,3).Value = 'Gombe'
$ws.Cells.Item(187,4).Value = 'Kumo'
$ws.Cells.Item(187,5).Value = 'Groundnut Gargaja'
$ws.Cells.Item(187,6).Value = 100000
$ws.Cells.Item(187,7).Value = 85
$ws.Cells.Item(187,8).Value = 1176.4705882352901
$ws.Cells.Item(187,9).Value = 'medium'
$ws.Cells.Item(187,10).Value = 'New'
$ws.Cells.Item(187,11).Value = 1179.4705882352901

$ws.Cells.Item(188,1).Value = 46051.749662708331
$ws.Cells.Item(188,2).Value = 'MH001GMB'
$ws.Cells.Item(188,3).Value = 'Gombe'
$ws.Cells.Item(188,4).Value = 'Kumo'
$ws.Cells.Item(188,5).Value = 'Groundnut Kampala'
$ws.Cells.Item(188,6).Value = 96000
$ws.Cells.Item(188,7).Value = 85
$ws.Cells.Item(188,8).Value = 1129.4117647058799
$ws.Cells.Item(188,9).Value = 'low'
$ws.Cells.Item(188,10).Value = 'New'
$ws.Cells.Item(188,11).Value = 1132.4117647058799

$ws.Cells.Item(189,1).Value = 46051.750264953713
$ws.Cells.Item(189,2).Value = 'MH001GMB'
$ws.Cells.Item(189,3).Value = 'Gombe'
$ws.Cells.Item(189,4).Value = 'Kumo'
$ws.Cells.Item(189,5).Value = 'Maize White'
$ws.Cells.Item(189,6).Value = 22000
$ws.Cells.Item(189,7).Value = 95
$ws.Cells.Item(189,8).Value = 231.57894736842101
$ws.Cells.Item(189,9).Value = 'high'
$ws.Cells.Item(189,10).Value = 'New'
$ws.Cells.Item(189,11).Value = 234.57894736842101

$ws.Cells.Item(190,1).Value = 46051.751193750002
$ws.Cells.Item(190,2).Value = 'MH001GMB'
$ws.Cells.Item(190,3).Value = 'Gombe'
$ws.Cells.Item(190,4).Value = 'Kumo'
$ws.Cells.Item(190,5).Value = 'Millet'
$ws.Cells.Item(190,6).Value = 30000
$ws.Cells.Item(190,7).Value = 100
$ws.Cells.Item(190,8).Value = 300
$ws.Cells.Item(190,9).Value = 'medium'
$ws.Cells.Item(190,10).Value = 'New'
$ws.Cells.Item(190,11).Value = 303

$ws.Cells.Item(191,1).Value = 46051.752038564817
$ws.Cells.Item(191,2).Value = 'MH001GMB'
$ws.Cells.Item(191,3).Value = 'Gombe'
$ws.Cells.Item(191,4).Value = 'Kumo'
$ws.Cells.Item(191,5).Value = 'Rice Paddy'
$ws.Cells.Item(191,6).Value = 28000
$ws.Cells.Item(191,7).Value = 75
$ws.Cells.Item(191,8).Value = 373.33333333333297
$ws.Cells.Item(191,9).Value = 'medium'
$ws.Cells.Item(191,10).Value = 'New'
$ws.Cells.Item(191,11).Value = 376.33333333333297

$ws.Cells.Item(192,1).Value = 46051.755584768522
$ws.Cells.Item(192,2).Value = 'MH001GMB'
$ws.Cells.Item(192,3).Value = 'Gombe'
$ws.Cells.Item(192,4).Value = 'Kumo'
$ws.Cells.Item(192,5).Value = 'Sorghum Red'
$ws.Cells.Item(192,6).Value = 23000
$ws.Cells.Item(192,7).Value = 100
$ws.Cells.Item(192,8).Value = 230
$ws.Cells.Item(192,9).Value = 'low'
$ws.Cells.Item(192,10).Value = 'New'
$ws.Cells.Item(192,11).Value = 233

$ws.Cells.Item(193,1).Value = 46051.756378518519
$ws.Cells.Item(193,2).Value = 'MH001GMB'
$ws.Cells.Item(193,3).Value = 'Gombe'
$ws.Cells.Item(193,4).Value = 'Kumo'
$ws.Cells.Item(193,5).Value = 'Sorghum White'
$ws.Cells.Item(193,6).Value = 23000
$ws.Cells.Item(193,7).Value = 100
$ws.Cells.Item(193,8).Value = 230
$ws.Cells.Item(193,9).Value = 'low'
$ws.Cells.Item(193,10).Value = 'New'
$ws.Cells.Item(193,11).Value = 233

$ws.Cells.Item(194,1).Value = 46051.756998287026
$ws.Cells.Item(194,2).Value = 'MH001GMB'
$ws.Cells.Item(194,3).Value = 'Gombe'
$ws.Cells.Item(194,4).Value = 'Kumo'
$ws.Cells.Item(194,5).Value = 'Soya Beans'
$ws.Cells.Item(194,6).Value = 62000
$ws.Cells.Item(194,7).Value = 95
$ws.Cells.Item(194,8).Value = 652.63157894736798
$ws.Cells.Item(194,9).Value = 'medium'
$ws.Cells.Item(194,10).Value = 'New'
$ws.Cells.Item(194,11).Value = 655.63157894736798

$ws.Cells.Item(195,1).Value = 46052.306498148151
$ws.Cells.Item(195,2).Value = 'MA001BOR'
$ws.Cells.Item(195,3).Value = 'Borno'
$ws.Cells.Item(195,4).Value = 'Lashe Money'
$ws.Cells.Item(195,5).Value = 'Cowpea White'
$ws.Cells.Item(195,6).Value = 68000
$ws.Cells.Item(195,7).Value = 105
$ws.Cells.Item(195,8).Value = 647.61904761904702
$ws.Cells.Item(195,9).Value = 'high'
$ws.Cells.Item(195,10).Value = 'New'
$ws.Cells.Item(195,11).Value = 650.61904761904702

$ws.Cells.Item(196,1).Value = 46052.30701689815
$ws.Cells.Item(196,2).Value = 'MA001BOR'
$ws.Cells.Item(196,3).Value = 'Borno'
$ws.Cells.Item(196,4).Value = 'Lashe Money'
$ws.Cells.Item(196,5).Value = 'Cowpea Brown'
$ws.Cells.Item(196,6).Value = 61000
$ws.Cells.Item(196,7).Value = 105
$ws.Cells.Item(196,8).Value = 580.95238095238096
$ws.Cells.Item(196,9).Value = 'high'
$ws.Cells.Item(196,10).Value = 'New'
$ws.Cells.Item(196,11).Value = 583.95238095238096

$ws.Cells.Item(197,1).Value = 46052.307482199067
$ws.Cells.Item(197,2).Value = 'MA001BOR'
$ws.Cells.Item(197,3).Value = 'Borno'
$ws.Cells.Item(197,4).Value = 'Lashe Money'
$ws.Cells.Item(197,5).Value = 'Rice Paddy'
$ws.Cells.Item(197,6).Value = 36000
$ws.Cells.Item(197,7).Value = 95
$ws.Cells.Item(197,8).Value = 378.94736842105198
$ws.Cells.Item(197,9).Value = 'high'
$ws.Cells.Item(197,10).Value = 'New'
$ws.Cells.Item(197,11).Value = 381.94736842105198

$ws.Cells.Item(198,1).Value = 46052.307945995373
$ws.Cells.Item(198,2).Value = 'MA001BOR'
$ws.Cells.Item(198,3).Value = 'Borno'
$ws.Cells.Item(198,4).Value = 'Lashe Money'
$ws.Cells.Item(198,5).Value = 'honeybeans'
$ws.Cells.Item(198,6).Value = 67000
$ws.Cells.Item(198,7).Value = 103
$ws.Cells.Item(198,8).Value = 650.48543689320297
$ws.Cells.Item(198,9).Value = 'high'
$ws.Cells.Item(198,10).Value = 'New'
$ws.Cells.Item(198,11).Value = 653.48543689320297

$ws.Cells.Item(199,1).Value = 46052.308440162044
$ws.Cells.Item(199,2).Value = 'MA001BOR'
$ws.Cells.Item(199,3).Value = 'Borno'
$ws.Cells.Item(199,4).Value = 'Lashe Money'
$ws.Cells.Item(199,5).Value = 'Soya Beans'
$ws.Cells.Item(199,6).Value = 71000
$ws.Cells.Item(199,7).Value = 115
$ws.Cells.Item(199,8).Value = 617.39130434782601
$ws.Cells.Item(199,9).Value = 'high'
$ws.Cells.Item(199,10).Value = 'New'
$ws.Cells.Item(199,11).Value = 620.39130434782601

$ws.Cells.Item(200,1).Value = 46051.797552916672
$ws.Cells.Item(200,2).Value = 'IS001KDN'
$ws.Cells.Item(200,3).Value = 'Kaduna'
$ws.Cells.Item(200,4).Value = 'Pambegua'
$ws.Cells.Item(200,5).Value = 'Soya Beans'
$ws.Cells.Item(200,6).Value = 60000
$ws.Cells.Item(200,7).Value = 100
$ws.Cells.Item(200,8).Value = 630
$ws.Cells.Item(200,9).Value = 'medium'
$ws.Cells.Item(200,10).Value = 'New'
$ws.Cells.Item(200,11).Value = 98.238095238095198

$ws.Cells.Item(201,1).Value = 46051.798799305558
$ws.Cells.Item(201,2).Value = 'IS001KDN'
$ws.Cells.Item(201,3).Value = 'Kaduna'
$ws.Cells.Item(201,4).Value = 'Pambegua'
$ws.Cells.Item(201,5).Value = 'Sorghum White'
$ws.Cells.Item(201,6).Value = 34000
$ws.Cells.Item(201,7).Value = 100
$ws.Cells.Item(201,8).Value = 260
$ws.Cells.Item(201,9).Value = 'medium'
$ws.Cells.Item(201,10).Value = 'New'
$ws.Cells.Item(201,11).Value = 133.76923076923001

$ws.Cells.Item(202,1).Value = 46051.80110222222
$ws.Cells.Item(202,2).Value = 'IS001KDN'
$ws.Cells.Item(202,3).Value = 'Kaduna'
$ws.Cells.Item(202,4).Value = 'Pambegua'
$ws.Cells.Item(202,5).Value = 'Maize'
$ws.Cells.Item(202,6).Value = 22000
$ws.Cells.Item(202,7).Value = 100
$ws.Cells.Item(202,8).Value = 230
$ws.Cells.Item(202,9).Value = 'high'
$ws.Cells.Item(202,10).Value = 'New'
$ws.Cells.Item(202,11).Value = 98.652173913043399

$ws.Cells.Item(203,1).Value = 46051.802656782413
$ws.Cells.Item(203,2).Value = 'IS001KDN'
$ws.Cells.Item(203,3).Value = 'Kaduna'
$ws.Cells.Item(203,4).Value = 'Pambegua'
$ws.Cells.Item(203,5).Value = 'Rice Paddy'
$ws.Cells.Item(203,6).Value = 34000
$ws.Cells.Item(203,7).Value = 80
$ws.Cells.Item(203,8).Value = 360
$ws.Cells.Item(203,9).Value = 'medium'
$ws.Cells.Item(203,10).Value = 'New'
$ws.Cells.Item(203,11).Value = 97.4444444444444


$ws.Range("L186").Select()
